# Update dSF (column F) values per repull of data / recomputed mean.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    4  = -4
    6  = -4
    7  = -1
    8  = -2
    10 = 8
    11 = -5
    12 = 2
    13 = 6
    16 = -1
    18 = 1
    19 = -1
    23 = -2
    25 = 1
    26 = 1
    27 = 3
    28 = -4
    29 = -2
    34 = -4
    36 = -3
    37 = -1
    39 = -1
    42 = 5
    43 = -2
    46 = -5
    50 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
